$d = $word.ActiveDocument

# --- 1) ANALYZER paragraph: highlight it green and fix "python?" -> "python" ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ANALYZER: python?`r") {
        # Fix the trailing "?" in the second run only (scope the Find to this
        # paragraph's Range so the bold "ANALYZER" run is left untouched).
        $pr = $p.Range
        $pr.Find.ClearFormatting()
        $pr.Find.Execute(": python?", $false, $false, $false, $false, $false, `
                          $true, 1, $false, ": python", 2)

        # Highlight the whole paragraph (including the paragraph-mark run
        # properties in <w:pPr><w:rPr>) bright green - wdBrightGreen = 4.
        $p.Range.Font.HighlightColorIndex = 4
        break
    }
}

# --- 2) PLANNER paragraph: append a new run after the ": python?" run ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "PLANNER: python?`r") {
        $r = $p.Range
        $r.Find.ClearFormatting()
        $r.Find.Execute(": python?", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData>' + `
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:body><w:p><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' + `
               '<w:t xml:space="preserve"> -&gt; gestione dei sintomi, scelta delle modalit&#224;</w:t>' + `
               '</w:r></w:p></w:body></w:document>' + `
               '</pkg:xmlData></pkg:part></pkg:package>'

        $r.InsertXML($xml)
        break
    }
}
